# The "codeforiati:group-name" (col D) and "codeforiati:group-code" (col E)
# columns were reordered upstream, so for every row (including the header)
# the value that used to be in column D now belongs in column E and vice
# versa. Swap the two columns' values in place, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    # .Value2 is used for reads (this host's .Value getter w/o args is unreliable)
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2

    $dCell.Value = $eVal
    $eCell.Value = $dVal
}
